$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''34.667.78'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '''  +0.84%  '
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = '''1.806.35'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '''  +0.16%  '
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = '''  -0.02%  '
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = '''225.19'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '''  -1.07%  '
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = '''0.604'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '''  +0.37%  '
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = '''  +0.01%  '
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = '''39.47'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '''  +8.69%  '
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = '''  -2.14%  '
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = '''0.0672'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '''  -2.77%  '
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = '''  +3.82%  '
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = '''2.066.89'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '''  +0.13%  '
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = '''1.815.85'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '''  +0.43%  '
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = '''10.95'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '''  -3.11%  '
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = '''0.636'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '''  -1.26%  '
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = '''34.652.26'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '''  +0.83%  '
$ws.Range("E16").Style = "Normal"
$ws.Range("E17").Value = '''  -1.55%  '
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = '''68.02'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '''  -2.99%  '
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = '''240.90'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '''  -1.44%  '
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = '''0.0₃0770'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '''  -1.94%  '
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = '''11.12'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '''  -2.84%  '
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = '''  -0.02%  '
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = '''  -1.84%  '
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = '''  -2.63%  '
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = '''171.71'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '''  +0.37%  '
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = '''7.68'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '''  -6.83%  '
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = '''17.50'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '''  +0.53%  '
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = '''  -0.34%  '
$ws.Range("E28").Style = "Normal"
$ws.Range("E29").Value = '''  -0.03%  '
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Value = '''  -1.39%  '
$ws.Range("E30").Style = "Normal"
$ws.Range("E31").Value = '''  -1.46%  '
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = '''  -1.35%  '
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = '''  -2.98%  '
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = '''  +1.05%  '
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Value = '''  -2.70%  '
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = '''1.308.86'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '''  -5.14%  '
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Value = '''  -0.25%  '
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = '''2.38'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '''  +1.66%  '
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Value = '''  +0.01%  '
$ws.Range("E39").Style = "Normal"
$ws.Range("B40").Value = '''InjectiveProtocol'
$ws.Range("B40").Style = "Normal"
$ws.Range("C40").Value = '''https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("C40").Style = "Normal"
$ws.Range("D40").Value = '''14.74'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '''  +10.45%  '
$ws.Range("E40").Style = "Normal"
$ws.Range("B41").Value = '''WEMIXToken'
$ws.Range("B41").Style = "Normal"
$ws.Range("C41").Value = '''https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("C41").Style = "Normal"
$ws.Range("D41").Value = '''1.24'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '''  +5.19%  '
$ws.Range("E41").Style = "Normal"
$ws.Range("B42").Value = '''Aave'
$ws.Range("B42").Style = "Normal"
$ws.Range("C42").Value = '''https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("C42").Style = "Normal"
$ws.Range("D42").Value = '''83.00'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '''  +0.98%  '
$ws.Range("E42").Style = "Normal"
$ws.Range("B43").Value = '''HuobiToken'
$ws.Range("B43").Style = "Normal"
$ws.Range("C43").Value = '''https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("C43").Style = "Normal"
$ws.Range("D43").Value = '''2.45'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '''  +1.13%  '
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = '''  +0.24%  '
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = '''0.944'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '''  -0.53%  '
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = '''0.0518'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '''  +3.81%  '
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = '''1.968.22'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '''  +0.12%  '
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = '''5.74'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '''  -3.64%  '
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = '''  -0.07%  '
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = '''101.85'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '''  -1.16%  '
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = '''  -0.40%  '
$ws.Range("E51").Style = "Normal"
